$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3, 24),
    @(6, 7, 18),
    @(7, 7, 8),
    @(24, 7, 47),
    @(28, 14, 20),
    @(26, 730, 15),
    @(27, 730, 17)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

$ws.Range("C10").Select() | Out-Null
